$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so numeric-looking strings are not
# auto-converted to Excel numbers (original data is stored as inline text strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.052.83'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '1.902.24'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '333.57'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4646'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").Value = '0.4135'
$ws.Range("E8").Value = '  +4.36%  '
$ws.Range("D9").Value = '47.68'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '0.07992'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").Value = '1.006'
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '1.905.96'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").Value = '5.943'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '7.089'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '89.20'
$ws.Range("E17").Value = '  -1.64%  '
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").Value = '0.06586'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '17.53'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '29.062.54'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '5.439'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").Value = '11.28'
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("D25").Value = '2.231'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").Value = '2.128.61'
$ws.Range("E26").Value = '  +2.01%  '
$ws.Range("D27").Value = '157.59'
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '2.117'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '5.433'
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").Value = '118.30'
$ws.Range("E31").Value = '  -1.57%  '
$ws.Range("D32").Value = '0.9821'
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").Value = '1.428'
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("D35").Value = '3.598'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("D36").Value = '5.303'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").Value = '0.06102'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '0.02246'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '8.367'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").Value = '1.173'
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '0.5800'
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("D43").Value = '10.18'
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").Value = '0.1823'
$ws.Range("E44").Value = '  -2.79%  '
$ws.Range("D45").Value = '1.264'
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("D46").Value = '2.318'
$ws.Range("E46").Value = '  +11.82%  '
$ws.Range("D47").Value = '0.5509'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("D48").Value = '12.07'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").Value = '1.914'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("D50").Value = '0.07046'
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("D51").Value = '47.02'
$ws.Range("E51").Value = '  +19.52%  '
